$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 308; existing rows 308:324 shift down to 309:325
$ws.Rows("308:308").Insert()

# Populate the newly inserted row 308 with the new weekly price-point record
$ws.Cells.Item(308, 1).Value  = 7
$ws.Cells.Item(308, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(308, 3).Value  = 'Ñuble'
$ws.Cells.Item(308, 4).Value  = 45041
$ws.Cells.Item(308, 5).Value  = 16
$ws.Cells.Item(308, 6).Value  = 100112043
$ws.Cells.Item(308, 7).Value  = 'Pepino ensalada'
$ws.Cells.Item(308, 8).Value  = 'Sin especificar'
$ws.Cells.Item(308, 9).Value  = 'Primera'
$ws.Cells.Item(308, 10).Value = 60
$ws.Cells.Item(308, 11).Value = 10000
$ws.Cells.Item(308, 12).Value = 10000
$ws.Cells.Item(308, 13).Value = 10000
$ws.Cells.Item(308, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(308, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(308, 16).Value = 167
$ws.Cells.Item(308, 17).Value = 60
$ws.Cells.Item(308, 18).Value = 'Hortaliza'
